# Refresh the cryptos list "Price" (D) and "Volume(1h)" (E) columns with
# the latest scraped figures, as produced by the scheduled GitHub Actions
# scraper run. All of these cells hold plain text (not numbers/percentages)
# so the exact source formatting (trailing zeros, the "  +x.xx%  " padding,
# thousands separators written with dots, etc.) must be preserved exactly.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Values that Excel would NOT mistake for a number (they already contain
# extra punctuation/spacing/subscripts) can be written directly as text.
$directUpdates = @(
    @{ Addr = 'D2'; Value = '41.077.97' },
    @{ Addr = 'E2'; Value = '  -2.11%  ' },
    @{ Addr = 'D3'; Value = '2.134.28' },
    @{ Addr = 'E3'; Value = '  -3.73%  ' },
    @{ Addr = 'E4'; Value = '  +0.02%  ' },
    @{ Addr = 'E5'; Value = '  -2.64%  ' },
    @{ Addr = 'E6'; Value = '  -4.44%  ' },
    @{ Addr = 'E7'; Value = '  -5.83%  ' },
    @{ Addr = 'E8'; Value = '  +0.10%  ' },
    @{ Addr = 'E9'; Value = '  -6.86%  ' },
    @{ Addr = 'E10'; Value = '  -8.99%  ' },
    @{ Addr = 'E11'; Value = '  -7.28%  ' },
    @{ Addr = 'E12'; Value = '  -7.24%  ' },
    @{ Addr = 'E13'; Value = '  -3.90%  ' },
    @{ Addr = 'E14'; Value = '  -6.38%  ' },
    @{ Addr = 'D15'; Value = '2.454.98' },
    @{ Addr = 'E15'; Value = '  -3.75%  ' },
    @{ Addr = 'E16'; Value = '  -0.71%  ' },
    @{ Addr = 'D17'; Value = '2.133.62' },
    @{ Addr = 'E17'; Value = '  -3.93%  ' },
    @{ Addr = 'E18'; Value = '  -7.57%  ' },
    @{ Addr = 'D19'; Value = '40.880.60' },
    @{ Addr = 'E19'; Value = '  -2.40%  ' },
    @{ Addr = 'D20'; Value = '0.0₃0990' },
    @{ Addr = 'E20'; Value = '  -7.52%  ' },
    @{ Addr = 'E21'; Value = '  -5.81%  ' },
    @{ Addr = 'E22'; Value = '  -8.45%  ' },
    @{ Addr = 'E23'; Value = '  -2.61%  ' },
    @{ Addr = 'E24'; Value = '  -12.63%  ' },
    @{ Addr = 'E25'; Value = '  -7.76%  ' },
    @{ Addr = 'E26'; Value = '  -0.17%  ' },
    @{ Addr = 'E27'; Value = '  -9.58%  ' },
    @{ Addr = 'E28'; Value = '  -12.62%  ' },
    @{ Addr = 'E29'; Value = '  -1.00%  ' },
    @{ Addr = 'E30'; Value = '  -6.04%  ' },
    @{ Addr = 'E31'; Value = '  +1.50%  ' },
    @{ Addr = 'E32'; Value = '  -5.04%  ' },
    @{ Addr = 'E33'; Value = '  +2.25%  ' },
    @{ Addr = 'E35'; Value = '  -10.83%  ' },
    @{ Addr = 'E36'; Value = '  -4.68%  ' },
    @{ Addr = 'E37'; Value = '  -4.96%  ' },
    @{ Addr = 'E38'; Value = '  -2.98%  ' },
    @{ Addr = 'E39'; Value = '  -4.75%  ' },
    @{ Addr = 'E40'; Value = '  -4.64%  ' },
    @{ Addr = 'E41'; Value = '  -16.86%  ' },
    @{ Addr = 'E42'; Value = '  -8.37%  ' },
    @{ Addr = 'E43'; Value = '  -12.35%  ' },
    @{ Addr = 'E44'; Value = '  -6.40%  ' },
    @{ Addr = 'E45'; Value = '  -6.70%  ' },
    @{ Addr = 'E46'; Value = '  -5.13%  ' },
    @{ Addr = 'E47'; Value = '  -8.23%  ' },
    @{ Addr = 'E48'; Value = '  -4.39%  ' },
    @{ Addr = 'E49'; Value = '  -6.11%  ' },
    @{ Addr = 'E50'; Value = '  -3.10%  ' },
    @{ Addr = 'E51'; Value = '  -10.77%  ' }
)

foreach ($u in $directUpdates) {
    $ws.Range($u.Addr).Value = $u.Value
}

# Values that look like plain numbers (e.g. "235.10", "69.02") must be
# forced into text storage first, otherwise Excel auto-converts them to
# numeric values and silently drops significant trailing zeros. Flip the
# cell to Text, write the literal string, then restore its original
# (default/General) style so the cell formatting is left unchanged.
$forcedTextUpdates = @(
    @{ Addr = 'D5'; Value = '235.10' },
    @{ Addr = 'D6'; Value = '0.598' },
    @{ Addr = 'D7'; Value = '69.02' },
    @{ Addr = 'D9'; Value = '0.564' },
    @{ Addr = 'D10'; Value = '38.53' },
    @{ Addr = 'D11'; Value = '0.0886' },
    @{ Addr = 'D12'; Value = '53.09' },
    @{ Addr = 'D13'; Value = '0.0991' },
    @{ Addr = 'D14'; Value = '6.55' },
    @{ Addr = 'D16'; Value = '14.15' },
    @{ Addr = 'D18'; Value = '0.772' },
    @{ Addr = 'D21'; Value = '68.51' },
    @{ Addr = 'D22'; Value = '5.66' },
    @{ Addr = 'D23'; Value = '223.73' },
    @{ Addr = 'D24'; Value = '9.41' },
    @{ Addr = 'D27'; Value = '10.48' },
    @{ Addr = 'D28'; Value = '3.27' },
    @{ Addr = 'D29'; Value = '2.16' },
    @{ Addr = 'D30'; Value = '2.14' },
    @{ Addr = 'D31'; Value = '170.02' },
    @{ Addr = 'D32'; Value = '19.46' },
    @{ Addr = 'D34'; Value = '0.0743' },
    @{ Addr = 'D35'; Value = '5.02' },
    @{ Addr = 'D38'; Value = '4.11' },
    @{ Addr = 'D39'; Value = '0.0286' },
    @{ Addr = 'D41'; Value = '11.56' },
    @{ Addr = 'D43'; Value = '56.96' },
    @{ Addr = 'D44'; Value = '0.185' },
    @{ Addr = 'D45'; Value = '8.14' },
    @{ Addr = 'D47'; Value = '96.21' },
    @{ Addr = 'D48'; Value = '1.07' },
    @{ Addr = 'D51'; Value = '2.12' }
)

foreach ($u in $forcedTextUpdates) {
    $cell = $ws.Range($u.Addr)
    $cell.NumberFormat = "@"
    $cell.Value = $u.Value
    $cell.Style = "Normal"
}
